$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds values that look numeric (e.g. "1.00",
# "0.999", "0.110") but must stay exactly as authored, including trailing
# zeros, matching the source inline-string cells. Marking those cells as
# Text before assignment stops Excel from re-interpreting/rounding them
# as numbers (e.g. "1.00" -> 1, "0.110" -> 0.11).
$ws.Range('D2').Value = '51.802.01'
$ws.Range('E2').Value = '  +2.31%  '
$ws.Range('D3').Value = '3.027.21'
$ws.Range('E3').Value = '  +3.90%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '380.63'
$ws.Range('E5').Value = '  +2.08%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '104.28'
$ws.Range('E6').Value = '  +5.30%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.548'
$ws.Range('E7').Value = '  +2.91%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.599'
$ws.Range('E9').Value = '  +4.52%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.94'
$ws.Range('E10').Value = '  +4.42%  '
$ws.Range('E11').Value = '  +0.19%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0862'
$ws.Range('E12').Value = '  +2.32%  '
$ws.Range('D13').Value = '3.483.53'
$ws.Range('E13').Value = '  +3.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.65'
$ws.Range('E14').Value = '  +4.58%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.84'
$ws.Range('E15').Value = '  +4.94%  '
$ws.Range('D16').Value = '3.035.57'
$ws.Range('E16').Value = '  +4.13%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '11.33'
$ws.Range('E17').Value = '  -5.22%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.997'
$ws.Range('E18').Value = '  +1.69%  '
$ws.Range('D19').Value = '51.815.26'
$ws.Range('E19').Value = '  +2.40%  '
$ws.Range('E20').Value = '  +3.60%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.62'
$ws.Range('E21').Value = '  +4.36%  '
$ws.Range('D22').Value = '0.0₃0966'
$ws.Range('E22').Value = '  +2.84%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.49'
$ws.Range('E23').Value = '  +2.33%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '269.93'
$ws.Range('E24').Value = '  +2.15%  '
$ws.Range('E25').Value = '  +4.08%  '
$ws.Range('E26').Value = '  +6.86%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.53'
$ws.Range('E27').Value = '  +6.47%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.171'
$ws.Range('E28').Value = '  +7.21%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '26.21'
$ws.Range('E30').Value = '  +3.90%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.110'
$ws.Range('E31').Value = '  +2.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '10.44'
$ws.Range('E32').Value = '  +6.42%  '
$ws.Range('B33').Value = 'InjectiveProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '34.76'
$ws.Range('E33').Value = '  +6.19%  '
$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '51.28'
$ws.Range('E34').Value = '  +2.03%  '
$ws.Range('E35').Value = '  +0.88%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0450'
$ws.Range('E36').Value = '  +5.92%  '
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.32'
$ws.Range('E38').Value = '  +9.42%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '17.34'
$ws.Range('E39').Value = '  +7.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.61'
$ws.Range('E40').Value = '  +8.93%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.86'
$ws.Range('E41').Value = '  +5.61%  '
$ws.Range('B42').Value = 'TheGraph'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.282'
$ws.Range('E42').Value = '  +10.30%  '
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.117'
$ws.Range('E43').Value = '  +2.46%  '
$ws.Range('B44').Value = 'Monero'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '127.99'
$ws.Range('E44').Value = '  +7.70%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.80'
$ws.Range('E45').Value = '  +14.41%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '21.69'
$ws.Range('E46').Value = '  +4.92%  '
$ws.Range('E47').Value = '  +1.33%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.37'
$ws.Range('E48').Value = '  +1.78%  '
$ws.Range('D49').Value = '2.047.02'
$ws.Range('E49').Value = '  +3.62%  '
$ws.Range('D50').Value = '3.320.61'
$ws.Range('E50').Value = '  +3.78%  '
$ws.Range('B51').Value = 'BEAM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0325'
$ws.Range('E51').Value = '  +3.71%  '
